$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the reference-number column (column B) with the new bracketed
# reference numbers from the latest manuscript revision. Column A (the
# author/year labels) and the header row are unchanged. Order matches the
# authoring order so newly-created shared strings land in the same spot.
$ws.Range("B3").Value = "Engel [76]"
$ws.Range("B4").Value = "Gordon et al. [45]"
$ws.Range("B5").Value = "Herrmann et al. [47]"
$ws.Range("B7").Value = "Juarez [78]"
$ws.Range("B8").Value = "Lustig [79]"
$ws.Range("B11").Value = "Kramer [46]"
$ws.Range("B6").Value = "Holobinko [77]"
$ws.Range("B10").Value = "Ueda & Bell [29]"
$ws.Range("B2").Value = "Chesson et al. [40]"
# B9 ("Regan [18]") is unchanged.

# Restore the active-cell selection to C7 (matches the saved sheet view).
$ws.Range("C7").Select()
